$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so values like "1.00" or "355.00"
# are not silently converted to numbers (losing trailing zeros / thousand dots).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '66.484.59'
$ws.Range('E2').Value = '  +4.33%  '

$ws.Range('D3').Value = '3.493.90'
$ws.Range('E3').Value = '  +2.24%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '591.08'
$ws.Range('E5').Value = '  +3.71%  '

$ws.Range('D6').Value = '168.64'
$ws.Range('E6').Value = '  +6.77%  '

$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').Value = '3.493.11'
$ws.Range('E8').Value = '  +2.16%  '

$ws.Range('D9').Value = '0.572'
$ws.Range('E9').Value = '  +3.22%  '

$ws.Range('E10').Value = '  +0.73%  '

$ws.Range('E11').Value = '  +5.13%  '

$ws.Range('D12').Value = '0.437'
$ws.Range('E12').Value = '  +3.55%  '

$ws.Range('D13').Value = '4.104.88'
$ws.Range('E13').Value = '  +2.40%  '

$ws.Range('E14').Value = '  +0.21%  '

$ws.Range('D15').Value = '28.05'
$ws.Range('E15').Value = '  +3.49%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.553.78'
$ws.Range('E16').Value = '  +4.31%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0000177'
$ws.Range('E17').Value = '  +2.85%  '

$ws.Range('D18').Value = '3.504.21'
$ws.Range('E18').Value = '  +2.00%  '

$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  +4.17%  '

$ws.Range('D20').Value = '13.99'
$ws.Range('E20').Value = '  +3.11%  '

$ws.Range('D21').Value = '388.67'
$ws.Range('E21').Value = '  +2.20%  '

$ws.Range('D22').Value = '7.94'
$ws.Range('E22').Value = '  +1.79%  '

$ws.Range('D23').Value = '73.44'
$ws.Range('E23').Value = '  +2.63%  '

$ws.Range('E24').Value = '  -0.08%  '

$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = '0.531'
$ws.Range('E25').Value = '  +3.40%  '

$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0000123'
$ws.Range('E26').Value = '  +6.45%  '

$ws.Range('D27').Value = '10.12'
$ws.Range('E27').Value = '  +5.37%  '

$ws.Range('E28').Value = '  +1.52%  '

$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('D30').Value = '6.37'
$ws.Range('E30').Value = '  +6.45%  '

$ws.Range('D31').Value = '1.46'
$ws.Range('E31').Value = '  +4.83%  '

$ws.Range('D32').Value = '2.05'
$ws.Range('E32').Value = '  +2.87%  '

$ws.Range('D33').Value = '23.52'
$ws.Range('E33').Value = '  +2.92%  '

$ws.Range('D34').Value = '7.43'
$ws.Range('E34').Value = '  +6.37%  '

$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('D36').Value = '1.61'
$ws.Range('E36').Value = '  +6.67%  '

$ws.Range('E37').Value = '  +0.48%  '

$ws.Range('D38').Value = '0.880'
$ws.Range('E38').Value = '  +4.00%  '

$ws.Range('D39').Value = '1.89'
$ws.Range('E39').Value = '  +3.87%  '

$ws.Range('D40').Value = '0.0750'
$ws.Range('E40').Value = '  +4.00%  '

$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '4.63'
$ws.Range('E41').Value = '  +5.63%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '26.38'
$ws.Range('E42').Value = '  +2.10%  '

$ws.Range('D43').Value = '2.813.05'
$ws.Range('E43').Value = '  +0.58%  '

$ws.Range('D44').Value = '6.61'
$ws.Range('E44').Value = '  +3.08%  '

$ws.Range('D45').Value = '43.27'
$ws.Range('E45').Value = '  +0.57%  '

$ws.Range('D46').Value = '26.33'
$ws.Range('E46').Value = '  -0.22%  '

$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = '355.00'
$ws.Range('E47').Value = '  +6.81%  '

$ws.Range('D48').Value = '2.52'
$ws.Range('E48').Value = '  +6.14%  '

$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0310'
$ws.Range('E49').Value = '  +3.37%  '

$ws.Range('D50').Value = '1.07'
$ws.Range('E50').Value = '  +2.93%  '

$ws.Range('D51').Value = '33.83'
$ws.Range('E51').Value = '  +14.78%  '
